# Update cryptos.xlsx "Sheet1" per the Dec 31 2022 GitHub Actions symbol-list refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($row, $col, $text)
    $cell = $ws.Cells.Item($row, $col)
    # Force Text number format so numeric-looking strings (e.g. "245.79",
    # "0.0002000") are kept verbatim instead of being coerced into floats.
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# --- Price-only refreshes (rows 2-9) ---------------------------------------
Set-TextCell 2 4 "245.79"
Set-TextCell 3 4 "26.13"
Set-TextCell 4 4 "5.082"
Set-TextCell 5 4 "0.05605"
Set-TextCell 6 4 "6.480"
Set-TextCell 7 4 "3.038"
Set-TextCell 8 4 "0.8134"
Set-TextCell 9 4 "0.8439"

# --- Rows 10-20: coin list reshuffled + prices updated ----------------------
$rows10to20 = @(
    @{ Row = 10; Coin = "One";                               Link = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one";                                         Price = "0.009755"; Vol = "9OneONEBestin24h" },
    @{ Row = 11; Coin = "WazirX";                             Link = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx";                                          Price = "0.1346";   Vol = "10WazirXWRX" },
    @{ Row = 12; Coin = "MandalaExchangeToken";                Link = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx";                      Price = "0.07008";  Vol = "11MandalaExchangeTokenMDX" },
    @{ Row = 13; Coin = "LiechtensteinCryptoassetsExchange";   Link = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx";              Price = "0.03216";  Vol = "12LiechtensteinCryptoassetsExchangeLCX" },
    @{ Row = 14; Coin = "BitrueCoin";                          Link = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr";                                    Price = "0.02796";  Vol = "13BitrueCoinBTR" },
    @{ Row = 15; Coin = "BitMartToken";                        Link = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx";                              Price = "0.09392";  Vol = "14BitMartTokenBMX" },
    @{ Row = 16; Coin = "BitForexToken";                       Link = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf";                              Price = "0.001525"; Vol = "15BitForexTokenBF" },
    @{ Row = 17; Coin = "TigerCash";                           Link = "https://coinranking.com/coin/6hIn06L2+tigercash-tch";                                      Price = "0.006138"; Vol = "16TigerCashTCH" },
    @{ Row = 18; Coin = "LEO";                                 Link = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo";                                        Price = "3.557";    Vol = "17LEOLEO" },
    @{ Row = 19; Coin = "BTSEToken";                           Link = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse";                                 Price = "2.118";    Vol = "18BTSETokenBTSE" },
    @{ Row = 20; Coin = "BitpandaEcosystemToken";              Link = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best";                        Price = "0.3188";   Vol = "19BitpandaEcosystemTokenBEST" }
)

foreach ($entry in $rows10to20) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.Coin
    $ws.Cells.Item($r, 3).Value = $entry.Link
    Set-TextCell $r 4 $entry.Price
    $ws.Cells.Item($r, 5).Value = $entry.Vol
}

# --- Scattered price/volume-label refreshes further down the sheet ---------
Set-TextCell 22 4 "3.756"
Set-TextCell 23 4 "0.04684"
Set-TextCell 25 4 "0.001247"
Set-TextCell 26 4 "0.004603"
$ws.Cells.Item(26, 5).Value = "25HotbitTokenHTB"
Set-TextCell 27 4 "0.00009603"

Set-TextCell 40 4 "0.03659"
Set-TextCell 41 4 "0.006152"
Set-TextCell 43 4 "0.002500"
Set-TextCell 44 4 "0.008762"
Set-TextCell 45 4 "0.00005292"
$ws.Cells.Item(47, 5).Value = "46CoinbaseStockTokenCOINWorstin24h"
Set-TextCell 48 4 "0.002058"
Set-TextCell 49 4 "0.00002100"
Set-TextCell 50 4 "0.0002000"
